$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new test-case block (4 rows) is appended at the bottom of the sheet,
# following the same layout used by every other block: a "requirement"
# row, a "preconditions" row, a "steps" row and an "expected result" row,
# separated from the previous block by one blank row (row 50 stays empty,
# just like rows 5, 10, 15, ... throughout the sheet).
$ws.Range("A51").Value = "Tesztelt követelmény"
$ws.Range("B51").Value = "SINGLE módban először a Primaryből bayből tüzelne, de hibás lesz a lövés"

$ws.Range("A52").Value = "Előfeltételek (kiindulási állapot)"
$ws.Range("B52").Value = "Mindkét torepdóba rakunk legalább 1-1 torpedót, de a második meghibásodása 100%"

$ws.Range("A53").Value = "A teszt lépései"
$ws.Range("B53").Value = "Kilővünk egy torpedot (pl. TORPEDO,SINGLE)"

$ws.Range("A54").Value = "Elvárt kimenet/eredmény"
$ws.Range("B54").Value = "SUCCESS a fenti command után"

# Match the saved view state: scrolled down towards the new rows, with the
# newly added B52 cell selected/active.
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B52").Select()
